$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

$rows = @(
    @("Equip-101", "PT-4", "PDT-9", 9, 9, 7, 5, 6),
    @("Equip-102", "PT-9", "PDT-5", 4, 9, 4, 7, 7),
    @("Equip-103", "PT-9", "PDT-6", 2, 6, 3, 4, 2),
    @("Equip-104", "PT-9", "PDT-2", 7, 4, 4, 6, 3),
    @("Equip-105", "PT-3", "PDT-8", 6, 8, 5, 9, 7),
    @("Equip-106", "PT-3", "PDT-9", 9, 2, 6, 9, 9),
    @("Equip-107", "PT-7", "PDT-4", 7, 5, 2, 2, 5),
    @("Equip-108", "PT-8", "PDT-9", 9, 2, 5, 8, 4)
)

$startRow = 10
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($j = 0; $j -lt $data.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $data[$j]
    }
}
